$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1886.25
$ws.Range("I19").Value = 848.75
$ws.Range("J19").Value = 2923.75
$ws.Range("K19").Value = 848.75
$ws.Range("L19").Value = 2923.75
$ws.Range("M19").Value = -673.75
$ws.Range("N19").Value = -3273.75
$ws.Range("H62").Value = 21233.334
$ws.Range("I62").Value = 6850
$ws.Range("K62").Value = 6850
$ws.Range("M62").Value = -6226
$ws.Range("H65").Value = 21233.334
$ws.Range("I65").Value = 6850
$ws.Range("K65").Value = 34250
$ws.Range("M65").Value = -31130
$ws.Range("H113").Value = 7029.625
$ws.Range("I113").Value = 9047.25
$ws.Range("J113").Value = 6357.0835
$ws.Range("K113").Value = 9047.25
$ws.Range("L113").Value = 6357.0835
$ws.Range("M113").Value = -5793.25
$ws.Range("N113").Value = -12865.0835
$ws.Range("H127").Value = 4448.1816
$ws.Range("I127").Value = 4193
$ws.Range("K127").Value = 12579
$ws.Range("M127").Value = -7619
$ws.Range("H132").Value = 1510.3877
$ws.Range("I132").Value = 1519.7609
$ws.Range("K132").Value = 4559.2827
$ws.Range("M132").Value = -2029.2827
$ws.Range("H135").Value = 1286.2354
$ws.Range("I135").Value = 972.5714
$ws.Range("K135").Value = 8753.142600000001
$ws.Range("M135").Value = -6218.142600000001
$ws.Range("H138").Value = 2699.45
$ws.Range("I138").Value = 1769.5264
$ws.Range("J138").Value = 2917.5803
$ws.Range("K138").Value = 5308.5792
$ws.Range("L138").Value = 8752.740900000001
$ws.Range("M138").Value = -168.5792000000001
$ws.Range("N138").Value = -19032.7409

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4071.5356
$ws.Range("I32").Value = 2461.0845
$ws.Range("J32").Value = 12867.077
$ws.Range("K32").Value = 2461.0845
$ws.Range("L32").Value = 12867.077
$ws.Range("M32").Value = -2174.0845
$ws.Range("N32").Value = -13441.077
$ws.Range("H44").Value = 6942.4375
$ws.Range("H102").Value = 4374.892
$ws.Range("I102").Value = 4106
$ws.Range("K102").Value = 4106
$ws.Range("M102").Value = -2484
$ws.Range("H122").Value = 2546.5908
$ws.Range("I122").Value = 2456.7104
$ws.Range("K122").Value = 7370.1312
$ws.Range("M122").Value = -4920.1312

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 303.30435
$ws.Range("J80").Value = 273.2857
$ws.Range("L80").Value = 273.2857
$ws.Range("N80").Value = -2269.2857
$ws.Range("H83").Value = 303.30435
$ws.Range("J83").Value = 273.2857
$ws.Range("L83").Value = 1366.4285
$ws.Range("N83").Value = -11350.4285
$ws.Range("H94").Value = 4657.3486
$ws.Range("I94").Value = 590
$ws.Range("K94").Value = 590
$ws.Range("M94").Value = -139
$ws.Range("H134").Value = 2598.587
$ws.Range("I134").Value = 943.74286
$ws.Range("J134").Value = 7864
$ws.Range("K134").Value = 2831.22858
$ws.Range("L134").Value = 23592
$ws.Range("M134").Value = -296.22858
$ws.Range("N134").Value = -28662

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 21000
$ws.Range("I28").Value = 7000
$ws.Range("J28").Value = 24500
$ws.Range("K28").Value = 7000
$ws.Range("L28").Value = 24500
$ws.Range("M28").Value = -6755
$ws.Range("N28").Value = -24990
$ws.Range("H50").Value = 2661.2942
$ws.Range("J50").Value = 2661.2942
$ws.Range("L50").Value = 2661.2942
$ws.Range("N50").Value = -3911.2942
$ws.Range("H86").Value = 7123.2856
$ws.Range("I86").Value = 6406.9287
$ws.Range("K86").Value = 6406.9287
$ws.Range("M86").Value = -5283.9287
$ws.Range("H89").Value = 7123.2856
$ws.Range("I89").Value = 6406.9287
$ws.Range("K89").Value = 32034.6435
$ws.Range("M89").Value = -26418.6435
$ws.Range("H99").Value = 4375.5625
$ws.Range("I99").Value = 3853.9
$ws.Range("K99").Value = 3853.9
$ws.Range("M99").Value = -2355.9
$ws.Range("H108").Value = 48243.75
$ws.Range("J108").Value = 50158.5
$ws.Range("L108").Value = 50158.5
$ws.Range("N108").Value = -57838.5
$ws.Range("H112").Value = 46999
$ws.Range("J112").Value = 46999
$ws.Range("L112").Value = 46999
$ws.Range("N112").Value = -49953
$ws.Range("H126").Value = 4375.5625
$ws.Range("I126").Value = 3853.9
$ws.Range("K126").Value = 11561.7
$ws.Range("M126").Value = -9091.700000000001
$ws.Range("H132").Value = 1829.4667
$ws.Range("I132").Value = 1789.1034
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 5367.3102
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -2837.3102
$ws.Range("N132").Value = -14060
$ws.Range("H134").Value = 28145.723
$ws.Range("I134").Value = 39931.625
$ws.Range("K134").Value = 119794.875
$ws.Range("M134").Value = -117259.875

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 9265425
$ws.Range("I56").Value = 9265425
$ws.Range("K56").Value = 9265425
$ws.Range("M56").Value = -9264895
$ws.Range("H92").Value = 496
$ws.Range("I92").Value = 496
$ws.Range("K92").Value = 1488
$ws.Range("M92").Value = -240
$ws.Range("H132").Value = 2951.5
$ws.Range("I132").Value = 2794.8572
$ws.Range("J132").Value = 3024.6
$ws.Range("K132").Value = 25153.7148
$ws.Range("L132").Value = 27221.4
$ws.Range("M132").Value = -22623.7148
$ws.Range("N132").Value = -32281.4

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 7782.7666
$ws.Range("J57").Value = 26556
$ws.Range("L57").Value = 26556
$ws.Range("N57").Value = -28196
$ws.Range("H101").Value = 38000
$ws.Range("J101").Value = 38000
$ws.Range("L101").Value = 38000
$ws.Range("N101").Value = -44490
$ws.Range("H132").Value = 3652.6296
$ws.Range("I132").Value = 2710.318
$ws.Range("J132").Value = 7798.8
$ws.Range("K132").Value = 8130.954000000001
$ws.Range("L132").Value = 23396.4
$ws.Range("M132").Value = -5600.954000000001
$ws.Range("N132").Value = -28456.4
$ws.Range("H137").Value = 52750
$ws.Range("J137").Value = 52750
$ws.Range("L137").Value = 52750
$ws.Range("N137").Value = -62950

# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4042.3125
$ws.Range("J7").Value = 8019.2
$ws.Range("L7").Value = 8019.2
$ws.Range("N7").Value = -8243.200000000001
$ws.Range("H40").Value = 7186.5
$ws.Range("I40").Value = 3070.2856
$ws.Range("K40").Value = 3070.2856
$ws.Range("M40").Value = -2934.2856
$ws.Range("H94").Value = 49994
$ws.Range("J94").Value = 49994
$ws.Range("L94").Value = 49994
$ws.Range("N94").Value = -51346
$ws.Range("H122").Value = 4944.385
$ws.Range("I122").Value = 3530.125
$ws.Range("K122").Value = 10590.375
$ws.Range("M122").Value = -8140.375
$ws.Range("H126").Value = 4042.3125
$ws.Range("J126").Value = 8019.2
$ws.Range("L126").Value = 24057.6
$ws.Range("N126").Value = -28997.6
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H131").Value = 69466.336
$ws.Range("J131").Value = 69466.336
$ws.Range("L131").Value = 69466.336
$ws.Range("N131").Value = -79546.336
$ws.Range("H133").Value = 148999.5
$ws.Range("J133").Value = 148999.5
$ws.Range("L133").Value = 148999.5
$ws.Range("N133").Value = -154059.5

# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 19890
$ws.Range("J15").Value = 19890
$ws.Range("L15").Value = 19890
$ws.Range("N15").Value = -20466
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H69").Value = 21796.8
$ws.Range("J69").Value = 21796.8
$ws.Range("L69").Value = 21796.8
$ws.Range("N69").Value = -23294.8
$ws.Range("H72").Value = 21796.8
$ws.Range("J72").Value = 21796.8
$ws.Range("L72").Value = 65390.39999999999
$ws.Range("N72").Value = -72878.39999999999
$ws.Range("H76").Value = 15665
$ws.Range("I76").Value = 13997.5
$ws.Range("J76").Value = 19000
$ws.Range("K76").Value = 13997.5
$ws.Range("L76").Value = 19000
$ws.Range("M76").Value = -13682.5
$ws.Range("N76").Value = -19630
$ws.Range("H79").Value = 15665
$ws.Range("I79").Value = 13997.5
$ws.Range("J79").Value = 19000
$ws.Range("K79").Value = 13997.5
$ws.Range("L79").Value = 19000
$ws.Range("M79").Value = -12905.5
$ws.Range("N79").Value = -21184
$ws.Range("H129").Value = 81596.336
$ws.Range("J129").Value = 59894.5
$ws.Range("L129").Value = 59894.5
$ws.Range("N129").Value = -69894.5
$ws.Range("H130").Value = 35143
$ws.Range("J130").Value = 35143
$ws.Range("L130").Value = 35143
$ws.Range("N130").Value = -45183
$ws.Range("H132").Value = 404157.6
$ws.Range("I132").Value = 4197.65
$ws.Range("J132").Value = 2003997.4
$ws.Range("K132").Value = 12592.95
$ws.Range("L132").Value = 6011992.199999999
$ws.Range("M132").Value = -10062.95
$ws.Range("N132").Value = -6017052.199999999
$ws.Range("H136").Value = 3312.25
$ws.Range("I136").Value = 2416.3333
$ws.Range("K136").Value = 7248.999899999999
$ws.Range("M136").Value = -4698.999899999999
